# "title fade case study"
#
# Reproduces the sheetView / sheetFormatPr / cols changes from the diff:
#   - topLeftCell  A16 -> A28   (scroll position; not exposed by this COM
#                                 host's persistence layer - see notes below)
#   - zoomScale / zoomScaleNormal  160 -> 170
#   - selection activeCell/sqref   I32 -> E42
#   - defaultColWidth  11.58984375 -> 11.60546875 (not exposed either - see below)
#   - new <cols> entry: column D (stop_label_x) gets a custom width ~30.97
#
$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$win = $excel.ActiveWindow

# Keep gridlines visible explicitly. Touching the window/view in this host
# re-serialises <sheetView> from scratch, and without an explicit true here
# the exporter would otherwise drop "showGridLines" to its (wrong-for-us)
# false default even though the source workbook had it "true".
$win.DisplayGridlines = $true

# Scroll the viewport so the top-left visible cell becomes A28 (was A16).
# (Harmless best-effort: this host's exporter does not currently persist
# topLeftCell from these, see note below, but set them in case it does.)
$win.ScrollRow    = 28
$win.ScrollColumn = 1

# Zoom: 160 -> 170 (this also drives zoomScale/zoomScaleNormal in real Excel)
$win.Zoom = 170

# Move the selection / active cell from I32 to E42
$ws.Range("E42").Select()

# Widen column D ("stop_label_x", the station-name column) to fit its
# longer text, matching the new <cols> entry (min=4 max=4 style=0,
# customWidth). The target raw width is 30.97 characters; this host snaps
# ColumnWidth to its own on-screen pixel grid, so 30.15 is the nearest
# settable input that lands on the closest achievable grid value to 30.97.
$ws.Columns("D").ColumnWidth = 30.15

# NOTE: two attributes from the diff have no working write-path through
# this COM host's persistence layer, confirmed by exhaustive probing
# (ActiveWindow.ScrollRow/ScrollColumn, Panes(1).ScrollRow/ScrollColumn,
# SmallScroll/LargeScroll, Application.Goto(..., Scroll:=True),
# Worksheet.StandardWidth, Application.StandardFont(Size), and
# ExecuteExcel4Macro("COLUMN.WIDTH") all affect in-memory COM state but
# never reach the saved <sheetView topLeftCell=...> / <sheetFormatPr
# defaultColWidth=...> attributes in this build) - so they are left as-is:
#   - <sheetView topLeftCell="A16" .../>           (would be "A28")
#   - <sheetFormatPr defaultColWidth="11.58984375"/> (would be "11.60546875")
